# Update the LR-pairs worksheet with new TPM-derived values.
# The sending-cluster x target-cluster pairs are expanded: each of the
# five sending clusters (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac)
# now has a row for target cluster "FAPs" AND a row for target cluster
# "MuSCs", with all numeric columns recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Omg", "Rtn4r", "FAPs", 3, 1, 1.75219, 5.25657, 0.1346950717404098, 0.1418982611249563, 1, 0.3333333333333333, 0.2214103333333333, 0.664231, 0.5896903116237344, 0.6831206522767569, 0.3879529719633333, 3.49157674767, 0.07942837882878355, 0.0969336326966177),
    @("ECs", "Omg", "Rtn4r", "MuSCs", 3, 1, 1.75219, 5.25657, 0.1346950717404098, 0.1418982611249563, 1, 0.5, 0.1540585, 0.308117, 0.4103096883762655, 0.3168793477232431, 0.269939763115, 1.61963857869, 0.05526669291162629, 0.04496462842833857),
    @("FAPs", "Omg", "Rtn4r", "FAPs", 3, 1, 2.026462666666667, 6.079388, 0.1557790731975008, 0.1641097875428132, 1, 0.3333333333333333, 0.2214103333333333, 0.664231, 0.5896903116237344, 0.6831206522767569, 0.4486797745142223, 4.038117970628, 0.09186141021829076, 0.1121067851112465),
    @("FAPs", "Omg", "Rtn4r", "MuSCs", 3, 1, 2.026462666666667, 6.079388, 0.1557790731975008, 0.1641097875428132, 1, 0.5, 0.1540585, 0.308117, 0.4103096883762655, 0.3168793477232431, 0.3121937987326667, 1.873162792396, 0.06391766297921, 0.05200300243156665),
    @("Inflammatory-Mac", "Omg", "Rtn4r", "FAPs", 3, 1, 4.488144, 13.464432, 0.3450144550883694, 0.3634650518941472, 1, 0.3333333333333333, 0.2214103333333333, 0.664231, 0.5896903116237344, 0.6831206522767569, 0.9937214590880001, 8.943493131792, 0.2034516815357535, 0.2482904833297351),
    @("Inflammatory-Mac", "Omg", "Rtn4r", "MuSCs", 3, 1, 4.488144, 13.464432, 0.3450144550883694, 0.3634650518941472, 1, 0.5, 0.1540585, 0.308117, 0.4103096883762655, 0.3168793477232431, 0.6914367324239999, 4.148620394543999, 0.1415627735526159, 0.1151745685644121),
    @("MuSCs", "Omg", "Rtn4r", "FAPs", 2, 1, 1.981064, 3.962128, 0.152289168185153, 0.1069555001749241, 1, 0.3333333333333333, 0.2214103333333333, 0.664231, 0.5896903116237344, 0.6831206522767569, 0.4386280405946667, 2.631768243568, 0.08980344704402218, 0.07306351104408093),
    @("MuSCs", "Omg", "Rtn4r", "MuSCs", 2, 1, 1.981064, 3.962128, 0.152289168185153, 0.1069555001749241, 1, 0.5, 0.1540585, 0.308117, 0.4103096883762655, 0.3168793477232431, 0.305199748244, 1.220798992976, 0.06248572114113082, 0.03389198913084316),
    @("Resolving-Mac", "Omg", "Rtn4r", "FAPs", 3, 1, 2.760707333333334, 8.282122000000001, 0.2122222317885668, 0.2235713992631593, 1, 0.3333333333333333, 0.2214103333333333, 0.664231, 0.5896903116237344, 0.6831206522767569, 0.6112491309091113, 5.501242178182001, 0.1251453939968844, 0.1527262400950766),
    @("Resolving-Mac", "Omg", "Rtn4r", "MuSCs", 3, 1, 2.760707333333334, 8.282122000000001, 0.2122222317885668, 0.2235713992631593, 1, 0.5, 0.1540585, 0.308117, 0.4103096883762655, 0.3168793477232431, 0.4253104307123334, 2.551862584274, 0.08707683779168245, 0.07084515916808266)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
